$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.225.63'
$ws.Range("E2").Value = '  +3.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.588.97'
$ws.Range("E3").Value = '  +3.37%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '626.22'
$ws.Range("E5").Value = '  +3.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.59'
$ws.Range("E6").Value = '  +6.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.595.85'
$ws.Range("E7").Value = '  +3.53%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  +8.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.45'
$ws.Range("E11").Value = '  +7.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  +5.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  +5.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.64'
$ws.Range("E14").Value = '  +8.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.193.86'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.600.64'
$ws.Range("E16").Value = '  +3.78%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.461.52'
$ws.Range("E17").Value = '  +4.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.117'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").Value = '  +5.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.10'
$ws.Range("E20").Value = '  +7.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.16'
$ws.Range("E21").Value = '  +13.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.08'
$ws.Range("E22").Value = '  +4.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.646'
$ws.Range("E23").Value = '  +3.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.66'
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000134'
$ws.Range("E25").Value = '  +8.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.72'
$ws.Range("E26").Value = '  +6.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.732.41'
$ws.Range("E27").Value = '  +3.44%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.32'
$ws.Range("E29").Value = '  +12.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.63'
$ws.Range("E30").Value = '  +4.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("E31").Value = '  +9.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.172'
$ws.Range("E32").Value = '  +7.70%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.53'
$ws.Range("E34").Value = '  +6.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.95'
$ws.Range("E35").Value = '  +5.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.47'
$ws.Range("E36").Value = '  +3.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.583.53'
$ws.Range("E37").Value = '  +3.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.43'
$ws.Range("E38").Value = '  +6.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.41'
$ws.Range("E39").Value = '  +11.00%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0925'
$ws.Range("E41").Value = '  +6.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '178.58'
$ws.Range("E42").Value = '  +3.77%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.70'
$ws.Range("E44").Value = '  +5.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.75'
$ws.Range("E45").Value = '  +22.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.910'
$ws.Range("E46").Value = '  +3.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.38'
$ws.Range("E47").Value = '  +12.71%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.94'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.74'
$ws.Range("E49").Value = '  +10.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.84'
$ws.Range("E50").Value = '  +3.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.268'
$ws.Range("E51").Value = '  +9.16%  '
